# Inserts a new data row just above the current row 92 ("Haba" price record),
# shifting the existing rows 92-169 down to 93-170, and fills the new row
# with its own values while leaving the rest of the rows (and their
# contents) otherwise untouched - Excel re-numbers the formerly-92..169
# rows to 93..170 automatically as part of the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above row 92 - this shifts rows 92:169 down to 93:170.
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new record's values.
$ws.Cells.Item(92, 1).Value = 4
$ws.Cells.Item(92, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(92, 3).Value = "Los Lagos"
$ws.Cells.Item(92, 4).Value = 45216
$ws.Cells.Item(92, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(92, 5).Value = 10
$ws.Cells.Item(92, 6).Value = 100112026
$ws.Cells.Item(92, 7).Value = "Haba"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 150
$ws.Cells.Item(92, 11).Value = 16000
$ws.Cells.Item(92, 12).Value = 16000
$ws.Cells.Item(92, 13).Value = 16000
$ws.Cells.Item(92, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(92, 15).Value = "Región Metropolitana"
$ws.Cells.Item(92, 16).Value = 640
$ws.Cells.Item(92, 17).Value = 25
$ws.Cells.Item(92, 18).Value = "Hortaliza"
